# Add the ability to import username and password for contacts.
#
# The "Data" worksheet gets two new trailing columns, K ("Username") and
# L ("Password"), each with an explanatory cell comment -- mirroring the
# existing "Email Message Categories" column (J) and its comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New header cells.
$ws.Range("K1").Value = "Username"
$ws.Range("L1").Value = "Password"

# Give the new headers the same look as the existing "required column"
# headers (J1 uses that style already) without minting a brand-new style.
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Explanatory comments, matching the style/author of the existing J1 comment.
$ws.Range("K1").AddComment("If this is present, password must also be present. Minimum length is 5 character.")
$ws.Range("L1").AddComment("If this is present, username must also be present. Minimum length is 8 characters.")

# Match the workbook's new active selection (L1 on the Data sheet).
$null = $ws.Activate()
$null = $ws.Range("L1").Select()
